# Auto-applied numeric corrections to market-value columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2967.3333
$ws.Range("J32").Value = 2967.3333
$ws.Range("L32").Value = 2967.3333
$ws.Range("N32").Value = -3619.3333

$ws.Range("H40").Value = 1930.6154
$ws.Range("J40").Value = 2014
$ws.Range("L40").Value = 2014
$ws.Range("N40").Value = -2364

$ws.Range("H132").Value = 4699458
$ws.Range("I132").Value = 5955314
$ws.Range("J132").Value = 10928
$ws.Range("K132").Value = 17865942
$ws.Range("L132").Value = 32784
$ws.Range("M132").Value = -17863412
$ws.Range("N132").Value = -37844

$ws.Range("H135").Value = 1216.3572
$ws.Range("I135").Value = 509.3
$ws.Range("K135").Value = 4583.7
$ws.Range("M135").Value = -2048.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7294.8823
$ws.Range("I2").Value = 1319.091
$ws.Range("K2").Value = 1319.091
$ws.Range("M2").Value = -1206.091

$ws.Range("H32").Value = 5428.909
$ws.Range("I32").Value = 5228.244
$ws.Range("J32").Value = 8171.3335
$ws.Range("K32").Value = 5228.244
$ws.Range("L32").Value = 8171.3335
$ws.Range("M32").Value = -4941.244
$ws.Range("N32").Value = -8745.333500000001

$ws.Range("H61").Value = 1208.8667
$ws.Range("I61").Value = 1110.5588
$ws.Range("J61").Value = 1512.7273
$ws.Range("K61").Value = 1110.5588
$ws.Range("L61").Value = 1512.7273
$ws.Range("M61").Value = -898.5588
$ws.Range("N61").Value = -1936.7273

$ws.Range("H74").Value = 1674.3334
$ws.Range("I74").Value = 785.8889
$ws.Range("K74").Value = 785.8889
$ws.Range("M74").Value = 88.11109999999996

$ws.Range("H77").Value = 1674.3334
$ws.Range("I77").Value = 785.8889
$ws.Range("K77").Value = 3929.4445
$ws.Range("M77").Value = 438.5554999999999

$ws.Range("H116").Value = 7294.8823
$ws.Range("I116").Value = 1319.091
$ws.Range("K116").Value = 1319.091
$ws.Range("M116").Value = 974.9090000000001

$ws.Range("H132").Value = 2074.6765
$ws.Range("I132").Value = 1881.7
$ws.Range("J132").Value = 2350.3572
$ws.Range("K132").Value = 5645.1
$ws.Range("L132").Value = 7051.071599999999
$ws.Range("M132").Value = -3115.1
$ws.Range("N132").Value = -12111.0716

$ws.Range("H136").Value = 1208.8667
$ws.Range("I136").Value = 1110.5588
$ws.Range("J136").Value = 1512.7273
$ws.Range("K136").Value = 3331.6764
$ws.Range("L136").Value = 4538.1819
$ws.Range("M136").Value = -781.6764000000003
$ws.Range("N136").Value = -9638.1819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7294.8823
$ws.Range("I3").Value = 1319.091
$ws.Range("K3").Value = 1319.091
$ws.Range("M3").Value = -1205.091

$ws.Range("H105").Value = 100002260
$ws.Range("I105").Value = 100002260
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 100002260
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -100000513
$ws.Range("N105").ClearContents()

$ws.Range("H134").Value = 4193.4883
$ws.Range("I134").Value = 1126.7059
$ws.Range("K134").Value = 3380.1177
$ws.Range("M134").Value = -845.1176999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 631.63635
$ws.Range("I58").Value = 543.13513
$ws.Range("J58").Value = 1099.4286
$ws.Range("K58").Value = 543.13513
$ws.Range("L58").Value = 1099.4286
$ws.Range("M58").Value = -340.13513
$ws.Range("N58").Value = -1505.4286

$ws.Range("H132").Value = 5609.9644
$ws.Range("I132").Value = 6103.6816
$ws.Range("J132").Value = 3799.6667
$ws.Range("K132").Value = 18311.0448
$ws.Range("L132").Value = 11399.0001
$ws.Range("M132").Value = -15781.0448
$ws.Range("N132").Value = -16459.0001

$ws.Range("H134").Value = 1011.5833
$ws.Range("I134").Value = 998.5122
$ws.Range("J134").Value = 1088.1428
$ws.Range("K134").Value = 2995.5366
$ws.Range("L134").Value = 3264.4284
$ws.Range("M134").Value = -460.5365999999999
$ws.Range("N134").Value = -8334.428400000001

$ws.Range("H136").Value = 631.63635
$ws.Range("I136").Value = 543.13513
$ws.Range("J136").Value = 1099.4286
$ws.Range("K136").Value = 1629.40539
$ws.Range("L136").Value = 3298.2858
$ws.Range("M136").Value = 920.5946100000001
$ws.Range("N136").Value = -8398.2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 909.25
$ws.Range("I122").Value = 640.5714
$ws.Range("J122").Value = 1019.8823
$ws.Range("K122").Value = 5765.1426
$ws.Range("L122").Value = 9178.940699999999
$ws.Range("M122").Value = -3315.1426
$ws.Range("N122").Value = -14078.9407

$ws.Range("H131").Value = 24391596
$ws.Range("I131").Value = 100000330
$ws.Range("J131").Value = 1683.3226
$ws.Range("K131").Value = 300000990
$ws.Range("L131").Value = 5049.9678
$ws.Range("M131").Value = -299995950
$ws.Range("N131").Value = -15129.9678

$ws.Range("H140").Value = 22666.2
$ws.Range("I140").Value = 57173.95
$ws.Range("K140").Value = 171521.85
$ws.Range("M140").Value = -166341.85

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 687.375
$ws.Range("I97").Value = 687.375
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 687.375
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -191.375
$ws.Range("N97").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2015.5714
$ws.Range("I7").Value = 1568
$ws.Range("J7").Value = 2351.25
$ws.Range("K7").Value = 1568
$ws.Range("L7").Value = 2351.25
$ws.Range("M7").Value = -1456
$ws.Range("N7").Value = -2575.25

$ws.Range("H46").Value = 6187.5
$ws.Range("J46").Value = 7916.6665
$ws.Range("L46").Value = 7916.6665
$ws.Range("N46").Value = -8292.666499999999

$ws.Range("H122").Value = 22739048
$ws.Range("I122").Value = 41685520
$ws.Range("J122").Value = 3280
$ws.Range("K122").Value = 125056560
$ws.Range("L122").Value = 9840
$ws.Range("M122").Value = -125054110
$ws.Range("N122").Value = -14740

$ws.Range("H126").Value = 2015.5714
$ws.Range("I126").Value = 1568
$ws.Range("J126").Value = 2351.25
$ws.Range("K126").Value = 4704
$ws.Range("L126").Value = 7053.75
$ws.Range("M126").Value = -2234
$ws.Range("N126").Value = -11993.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 424.3409
$ws.Range("I136").Value = 244.21428
$ws.Range("J136").Value = 739.5625
$ws.Range("K136").Value = 732.64284
$ws.Range("L136").Value = 2218.6875
$ws.Range("M136").Value = 1817.35716
$ws.Range("N136").Value = -7318.6875
